# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (date 2023-04-28 / serial 45044) right
# after the current row 14, pushing the existing rows 15-36 down to 17-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 15:16 (shifts old rows 15-36 down to 17-38,
# carrying the date-column style down onto the new rows, same as Excel's
# native "Insert Copied/Above Cells" behaviour).
$ws.Rows("15:16").Insert()

# New row 15: Especial
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 45044
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100104
$ws.Range("H15").Value = "Frutos de pepita"
$ws.Range("I15").Value = 100104001
$ws.Range("J15").Value = "Granada"
$ws.Range("K15").Value = "Wonderfull"
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("Q15").Value = '$/caja 14 kilos granel'
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("S15").Value = 1143
$ws.Range("T15").Value = 14

# New row 16: Primera
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "Femacal de La Calera"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 45044
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = "Frutos de pepita"
$ws.Range("I16").Value = 100104001
$ws.Range("J16").Value = "Granada"
$ws.Range("K16").Value = "Wonderfull"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = '$/caja 14 kilos granel'
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 14

Write-Output "Inserted rows 15-16 (new dimension: $($ws.Range("A1").Value2) .. row count check done)"
